$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.02
$data[1,0] = 1.02
$data[2,0] = 1.02
$data[3,0] = 1.02
$data[4,0] = 1.02
$data[5,0] = 1.02
$data[6,0] = 1.02
$data[7,0] = 1.02
$data[8,0] = 1.02
$data[9,0] = 1.02
$data[10,0] = 1.02
$data[11,0] = 1.02
$data[12,0] = 1.02
$data[13,0] = 1.02
$data[14,0] = 1.02
$data[15,0] = 1.02
$data[16,0] = 1.02
$data[17,0] = 1.02
$data[18,0] = 1.02
$data[19,0] = 1.02
$data[20,0] = 1.02
$data[21,0] = 1.02
$data[22,0] = 1.02
$data[23,0] = 1.02
$ws.Range("B2:B25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.010622732913008
$data[1,0] = 1.012267720130226
$data[2,0] = 1.013329534166126
$data[3,0] = 1.013775313124033
$data[4,0] = 1.013850126184366
$data[5,0] = 1.013335493053048
$data[6,0] = 1.011179212277531
$data[7,0] = 1.007358940637753
$data[8,0] = 1.004797244699636
$data[9,0] = 1.003684261183983
$data[10,0] = 1.00327026682997
$data[11,0] = 1.003359096700946
$data[12,0] = 1.003650052279565
$data[13,0] = 1.003829241825999
$data[14,0] = 1.004871029079135
$data[15,0] = 1.005523496783809
$data[16,0] = 1.00590370941189
$data[17,0] = 1.006033291474622
$data[18,0] = 1.005453530624436
$data[19,0] = 1.00356438935315
$data[20,0] = 1.002373230232305
$data[21,0] = 1.00300501310331
$data[22,0] = 1.005485146434493
$data[23,0] = 1.008349117105378
$ws.Range("C2:C25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.013268254309555
$data[1,0] = 1.014824216018873
$data[2,0] = 1.015828826175059
$data[3,0] = 1.016250651085927
$data[4,0] = 1.016321447647368
$data[5,0] = 1.015834464617906
$data[6,0] = 1.013794562926314
$data[7,0] = 1.010182517007348
$data[8,0] = 1.007761862039964
$data[9,0] = 1.006710501723597
$data[10,0] = 1.00631948184901
$data[11,0] = 1.00640337974255
$data[12,0] = 1.006678190153545
$data[13,0] = 1.006847443437462
$data[14,0] = 1.007831568482541
$data[15,0] = 1.00844801528546
$data[16,0] = 1.008807270608161
$data[17,0] = 1.008929715713359
$data[18,0] = 1.008381908323377
$data[19,0] = 1.006597279217181
$data[20,0] = 1.005472323144595
$data[21,0] = 1.0060689632019
$data[22,0] = 1.008411780164375
$data[23,0] = 1.011118477954687
$ws.Range("D2:D25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.012829630966766
$data[1,0] = 1.014248996179511
$data[2,0] = 1.015165307789247
$data[3,0] = 1.015550031501989
$data[4,0] = 1.015614599704553
$data[5,0] = 1.015170450409669
$data[6,0] = 1.013309756813113
$data[7,0] = 1.010014252119089
$data[8,0] = 1.007805237397172
$data[9,0] = 1.00684568862619
$data[10,0] = 1.006488799212478
$data[11,0] = 1.006565374725977
$data[12,0] = 1.006816197716495
$data[13,0] = 1.006970675272938
$data[14,0] = 1.00786885436316
$data[15,0] = 1.008431436462027
$data[16,0] = 1.008759289670581
$data[17,0] = 1.008871030162452
$data[18,0] = 1.008371106976771
$data[19,0] = 1.006742349700188
$data[20,0] = 1.005715556942551
$data[21,0] = 1.006260142868672
$data[22,0] = 1.008398368174521
$data[23,0] = 1.010868279412308
$ws.Range("E2:E25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.008866401775824
$data[1,0] = 1.011170520799082
$data[2,0] = 1.012658018727508
$data[3,0] = 1.013282571957314
$data[4,0] = 1.0133873916259
$data[5,0] = 1.01266636710805
$data[6,0] = 1.009645814975659
$data[7,0] = 1.004295744692175
$data[8,0] = 1.000708777428093
$data[9,0] = 0.9991503778886722
$data[10,0] = 0.998570700293726
$data[11,0] = 0.9986950807090049
$data[12,0] = 0.9991024784588267
$data[13,0] = 0.9993533801425385
$data[14,0] = 1.000812090378617
$data[15,0] = 1.001725679268545
$data[16,0] = 1.002258059003729
$data[17,0] = 1.002439502759521
$data[18,0] = 1.001627711917561
$data[19,0] = 0.9989825328530689
$data[20,0] = 0.9973146500667875
$data[21,0] = 0.9981992885366439
$data[22,0] = 1.001671980694062
$data[23,0] = 1.005682303827192
$ws.Range("F2:F25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.023594999628091
$data[1,0] = 1.023504579208683
$data[2,0] = 1.023443335729026
$data[3,0] = 1.023416932628374
$data[4,0] = 1.023412460900047
$data[5,0] = 1.02344298551034
$data[6,0] = 1.023565007130014
$data[7,0] = 1.023759139049229
$data[8,0] = 1.023874603872166
$data[9,0] = 1.023921309042974
$data[10,0] = 1.023938164268129
$data[11,0] = 1.023934571058996
$data[12,0] = 1.023922712353276
$data[13,0] = 1.023915340507297
$data[14,0] = 1.023871434968932
$data[15,0] = 1.023843013862359
$data[16,0] = 1.023826118601784
$data[17,0] = 1.023820303822167
$data[18,0] = 1.023846096030143
$data[19,0] = 1.02392621804881
$data[20,0] = 1.023973741382756
$data[21,0] = 1.023948818288664
$data[22,0] = 1.023844704313072
$data[23,0] = 1.02371141742819
$ws.Range("I2:I25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.015875912179003
$data[1,0] = 1.017150088860457
$data[2,0] = 1.017971744485485
$data[3,0] = 1.018316505173813
$data[4,0] = 1.018374353465014
$data[5,0] = 1.017976353784442
$data[6,0] = 1.016307117841986
$data[7,0] = 1.013343546487521
$data[8,0] = 1.011352153591792
$data[9,0] = 1.010485965543978
$data[10,0] = 1.010163623420673
$data[11,0] = 1.010232794355919
$data[12,0] = 1.010459332998333
$data[13,0] = 1.010598830796842
$data[14,0] = 1.011409556147989
$data[15,0] = 1.011917047514306
$data[16,0] = 1.012212683207891
$data[17,0] = 1.012313424113664
$data[18,0] = 1.011862637435466
$data[19,0] = 1.010392639736102
$data[20,0] = 1.009464907179711
$data[21,0] = 1.009957051354936
$data[22,0] = 1.011887224167813
$data[23,0] = 1.014112406671109
$ws.Range("J2:J25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.016129916656902
$data[1,0] = 1.017488479125081
$data[2,0] = 1.018364895301191
$data[3,0] = 1.018732715462583
$data[4,0] = 1.018794437791005
$data[5,0] = 1.018369812571446
$data[6,0] = 1.016589608609778
$data[7,0] = 1.013431680646141
$data[8,0] = 1.011311481969041
$data[9,0] = 1.010389694353909
$data[10,0] = 1.010046725662995
$data[11,0] = 1.010120319918557
$data[12,0] = 1.010361356291172
$data[13,0] = 1.010509789992768
$data[14,0] = 1.011372578105769
$data[15,0] = 1.011912773292779
$data[16,0] = 1.012227501417767
$data[17,0] = 1.012334755224157
$data[18,0] = 1.011854852665216
$data[19,0] = 1.01029039310712
$data[20,0] = 1.009303416377554
$data[21,0] = 1.009826953201412
$data[22,0] = 1.011881025611365
$data[23,0] = 1.014250649547849
$ws.Range("K2:K25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.015692611750708
$data[1,0] = 1.016914860441043
$data[2,0] = 1.017703128852677
$data[3,0] = 1.018033903713064
$data[4,0] = 1.01808940663465
$data[5,0] = 1.017707551076573
$data[6,0] = 1.016106223262073
$data[7,0] = 1.013263996046943
$data[8,0] = 1.011354692591863
$data[9,0] = 1.0105243472595
$data[10,0] = 1.010215364615088
$data[11,0] = 1.010281667800164
$data[12,0] = 1.010498818069175
$data[13,0] = 1.01063253749945
$data[14,0] = 1.011409722754907
$data[15,0] = 1.011896255789868
$data[16,0] = 1.012179695460119
$data[17,0] = 1.012276282654583
$data[18,0] = 1.011844091372952
$data[19,0] = 1.010434888172092
$data[20,0] = 1.009545645695088
$data[21,0] = 1.010017359913397
$data[22,0] = 1.011867663298968
$data[23,0] = 1.015552562608901
$ws.Range("L2:L25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.011741352102767
$data[1,0] = 1.013844987741668
$data[2,0] = 1.015202482377502
$data[3,0] = 1.015772313472753
$data[4,0] = 1.015867940953244
$data[5,0] = 1.015210099831563
$data[6,0] = 1.012453066933474
$data[7,0] = 1.00756533711616
$data[8,0] = 1.004285394863481
$data[9,0] = 1.002859689474238
$data[10,0] = 1.002329265808783
$data[11,0] = 1.002443082606589
$data[12,0] = 1.002815862078011
$data[13,0] = 1.003045429803249
$data[14,0] = 1.004379896328307
$data[15,0] = 1.005215485621801
$data[16,0] = 1.005702345219157
$data[17,0] = 1.005868263575753
$data[18,0] = 1.005125889380989
$data[19,0] = 1.002706111623407
$data[20,0] = 1.00117974896345
$data[21,0] = 1.001989382774208
$data[22,0] = 1.005166375688317
$data[23,0] = 1.008832593294169
$ws.Range("M2:M25").Value = $data

$data = New-Object "object[,]" 24,1
$data[0,0] = 1.017318572496894
$data[1,0] = 1.018594558655385
$data[2,0] = 1.019417381125646
$data[3,0] = 1.019762631413695
$data[4,0] = 1.019820561856106
$data[5,0] = 1.019421996970336
$data[6,0] = 1.017750390521378
$data[7,0] = 1.014782610555662
$data[8,0] = 1.012788389653606
$data[9,0] = 1.01192097151941
$data[10,0] = 1.011598171633318
$data[11,0] = 1.011667440799226
$data[12,0] = 1.011894301152496
$data[13,0] = 1.012033997053879
$data[14,0] = 1.012845873728017
$data[15,0] = 1.013354085790284
$data[16,0] = 1.013650141320464
$data[17,0] = 1.013751025289878
$data[18,0] = 1.013299598442891
$data[19,0] = 1.011827513178184
$data[20,0] = 1.01089846313515
$data[21,0] = 1.011391306211553
$data[22,0] = 1.013324220091218
$data[23,0] = 1.015552562608901
$ws.Range("N2:N25").Value = $data
